$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 5-7 (the MuSCs-as-sending-cluster block is gone in the
# new TPM run) before rewriting the remaining rows' numbers.
$ws.Range("A5:A7").EntireRow.Delete()

# Row 2: FAPs -> F12 -> Gp1ba -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "F12"
$ws.Range("C2").Value = "Gp1ba"
$ws.Range("D2").Value = "ECs"
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.8229573333333334
$ws.Range("N2").Value = 2.468872
$ws.Range("O2").Value = 0.2440777672676426
$ws.Range("P2").Value = 0.2440777672676426
$ws.Range("Q2").Value = 0.09943107660888889
$ws.Range("R2").Value = 0.8948796894800001
$ws.Range("S2").Value = 0.2440777672676426
$ws.Range("T2").Value = 0.2440777672676426

# Row 3: FAPs -> F12 -> Gp1ba -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "F12"
$ws.Range("C3").Value = "Gp1ba"
$ws.Range("D3").Value = "FAPs"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.4345811965947162
$ws.Range("P3").Value = 0.4345811965947162
$ws.Range("S3").Value = 0.4345811965947162
$ws.Range("T3").Value = 0.4345811965947162

# Row 4: FAPs -> F12 -> Gp1ba -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F12"
$ws.Range("C4").Value = "Gp1ba"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 1.083466
$ws.Range("N4").Value = 3.250398
$ws.Range("O4").Value = 0.3213410361376413
$ws.Range("P4").Value = 0.3213410361376413
$ws.Range("Q4").Value = 0.1309061678966666
$ws.Range("R4").Value = 1.17815551107
$ws.Range("S4").Value = 0.3213410361376413
$ws.Range("T4").Value = 0.3213410361376413
